# "aton done 3 problems"
# On the "Easy" sheet, Aton (PIC) finished 3 more Stack & Queue problems -
# rows 17-19 ("Implement Queue using Stacks", "Next Greater Element I",
# "Min Stack") - moving their Status from "On-going" to "Done". Copy the
# cell formatting (font/color) already used for "Done" elsewhere on the
# sheet (e.g. E2) onto those cells, then write the new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Easy")

$doneTemplate = $ws.Range("E2")
$targets = @("E17", "E18", "E19")

$doneTemplate.Copy()
foreach ($addr in $targets) {
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

foreach ($addr in $targets) {
    $ws.Range($addr).Value = "Done"
}

# Leave the selection where the user's last edit was.
[void]$ws.Range("E19").Select()
